$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "Ready for handoff" status text everywhere it appears to the
#    new handed-back status. We touch every cell that currently shows the
#    old text (Overview!B2:C3 plus B2/B3 on each language sheet) so they all
#    converge on a single shared string (the exporter compacts/removes the
#    now-unused old string automatically on save).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Per-language-sheet updates (zh-cn, de-de): mark status handed back,
#    populate "Latest Target File" (E) / "Latest Handback File" (F) with the
#    same file + hyperlink already used for the handoff columns (A / C), and
#    stamp the "Latest Handback DateTime" (G) column.
# ---------------------------------------------------------------------------
function Set-HandbackRow(
    $ws,
    [string]$targetRow,
    [string]$targetName,
    [string]$targetUrl,
    [string]$handbackName,
    [string]$handbackUrl,
    [string]$handbackDateTime
) {
    $eCell = "E" + $targetRow
    $fCell = "F" + $targetRow
    $gCell = "G" + $targetRow

    $ws.Hyperlinks.Add($ws.Range($eCell), $targetUrl, "", "", $targetName) | Out-Null
    $ws.Range($eCell).Font.Underline = 2
    $ws.Range($eCell).Font.Color = 15570276

    $ws.Hyperlinks.Add($ws.Range($fCell), $handbackUrl, "", "", $handbackName) | Out-Null
    $ws.Range($fCell).Font.Underline = 2
    $ws.Range($fCell).Font.Color = 15570276

    $ws.Range($gCell).Value = $handbackDateTime
}

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

Set-HandbackRow $wsZh "2" `
    "14308cb6-45cf-45fd-9604-96f88ed76f23.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/5caa1d525efc79e1126b35a53c2dcfa999c39aea/e2e/14308cb6-45cf-45fd-9604-96f88ed76f23.md" `
    "14308cb6-45cf-45fd-9604-96f88ed76f23.ec6f30ab863b55b04127d3e6c404cc409053bef1.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c732e24ac9c4e327bf63d8ded6295f2f114a34f9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/14308cb6-45cf-45fd-9604-96f88ed76f23.ec6f30ab863b55b04127d3e6c404cc409053bef1.zh-cn.xlf" `
    "2016-01-26 06:08:20"

Set-HandbackRow $wsZh "3" `
    "4b4d62aa-5dd3-435e-bfa5-547e514725f4.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/5caa1d525efc79e1126b35a53c2dcfa999c39aea/e2e/4b4d62aa-5dd3-435e-bfa5-547e514725f4.md" `
    "4b4d62aa-5dd3-435e-bfa5-547e514725f4.e7003ccfcd9e270c30a9325863f606f080f7548e.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c732e24ac9c4e327bf63d8ded6295f2f114a34f9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/4b4d62aa-5dd3-435e-bfa5-547e514725f4.e7003ccfcd9e270c30a9325863f606f080f7548e.zh-cn.xlf" `
    "2016-01-26 06:08:20"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

Set-HandbackRow $wsDe "2" `
    "14308cb6-45cf-45fd-9604-96f88ed76f23.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/5caa1d525efc79e1126b35a53c2dcfa999c39aea/e2e/14308cb6-45cf-45fd-9604-96f88ed76f23.md" `
    "14308cb6-45cf-45fd-9604-96f88ed76f23.ec6f30ab863b55b04127d3e6c404cc409053bef1.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f65afaf442afb712eb50a4f3ae19633b27cf6084/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/14308cb6-45cf-45fd-9604-96f88ed76f23.ec6f30ab863b55b04127d3e6c404cc409053bef1.de-de.xlf" `
    "2016-01-26 06:08:42"

Set-HandbackRow $wsDe "3" `
    "4b4d62aa-5dd3-435e-bfa5-547e514725f4.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/5caa1d525efc79e1126b35a53c2dcfa999c39aea/e2e/4b4d62aa-5dd3-435e-bfa5-547e514725f4.md" `
    "4b4d62aa-5dd3-435e-bfa5-547e514725f4.e7003ccfcd9e270c30a9325863f606f080f7548e.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f65afaf442afb712eb50a4f3ae19633b27cf6084/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/4b4d62aa-5dd3-435e-bfa5-547e514725f4.e7003ccfcd9e270c30a9325863f606f080f7548e.de-de.xlf" `
    "2016-01-26 06:08:42"

Write-Host "Handback report generated"
